# Auto update Excel log
# Appends the latest batch of sensor readings captured at ~2026-01-30 15:53-15:54
# to each per-sensor log sheet in the SeniorConnect master log workbook.

function Add-LogRow {
    param(
        $Worksheet,
        [string[]]$Values
    )

    $lastRow = $Worksheet.UsedRange.Rows.Count
    $newRow = $lastRow + 1

    for ($col = 1; $col -le $Values.Length; $col++) {
        $cell = $Worksheet.Cells.Item($newRow, $col)
        # Force text so dates/times stay as literal strings instead of
        # being re-interpreted as Excel date/time serials.
        $cell.NumberFormat = "@"
        $cell.Value = $Values[$col - 1]
    }
}

$wb = $excel.ActiveWorkbook

# mmWave presence sensor - one new reading
$wsMmWave = $wb.Worksheets.Item("mmWave")
Add-LogRow $wsMmWave @("2026-01-30", "15:53:49", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")

# PIR motion sensor - five new readings
$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRow $wsPIR @("2026-01-30", "15:53:49", "15:00", "Bathroom", "No Motion", "Inactive")
Add-LogRow $wsPIR @("2026-01-30", "15:53:52", "15:00", "Bathroom", "No Motion", "Inactive")
Add-LogRow $wsPIR @("2026-01-30", "15:53:57", "15:00", "Bathroom", "No Motion", "Inactive")
Add-LogRow $wsPIR @("2026-01-30", "15:54:02", "15:00", "Bathroom", "No Motion", "Inactive")
Add-LogRow $wsPIR @("2026-01-30", "15:54:07", "15:00", "Bathroom", "No Motion", "Inactive")

# Humidity sensor - three new readings
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRow $wsHumidity @("2026-01-30", "15:53:52", "15:00", "Bathroom", "87.7%", "Active")
Add-LogRow $wsHumidity @("2026-01-30", "15:53:57", "15:00", "Bathroom", "87.8%", "Active")
Add-LogRow $wsHumidity @("2026-01-30", "15:54:07", "15:00", "Bathroom", "87.8%", "Active")

# Proximity / door sensor - one new reading
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity @("2026-01-30", "15:53:54", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")

# Camera sensor - one new reading
$wsCamera = $wb.Worksheets.Item("Camera")
Add-LogRow $wsCamera @("2026-01-30", "15:53:54", "15:00", "Living Room Main Door", "Image Captured (ENTER)", "Active")
